$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row 27 with the new experiment entry
$ws.Range("A27").Value = 1.8
$ws.Range("B27").Value = 0.3
$ws.Range("C27").Value = 250
$ws.Range("D27").Value = 25
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = "nach ca 64 M Comp. => 1.0 rating "

# Match the numeric formatting style of column A (2 decimal places)
$ws.Range("A27").NumberFormat = "0.00"

# Update the active cell / selection to reflect the new last entry
$ws.Range("I27").Select()
